# Weekly update: insert a new price record for "Poroto verde" at
# Terminal La Palmera de La Serena, pushing the existing history down
# by one row (newest record goes on top, right under the header-adjacent
# block of rows that already sits at the top of the table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 58 - this shifts the existing rows 58:126 down
# to 59:127 (values, formatting and styles all move with the insert, the
# same way Excel's own Insert Row command behaves).
$ws.Rows.Item(58).Insert()

# Populate the newly-inserted (now blank) row 58 with the new weekly
# record. It mirrors the row that used to be at 58 (same market,
# region, product, variety, quality, unit, origin) except for the new
# sample date and the updated volume/price figures.
$ws.Range("A58").Value = 8
$ws.Range("B58").Value = "Terminal La Palmera de La Serena"
$ws.Range("C58").Value = "Coquimbo"
$ws.Range("D58").Value = 44482
$ws.Range("E58").Value = 4
$ws.Range("F58").Value = 100112031
$ws.Range("G58").Value = "Poroto verde"
$ws.Range("H58").Value = "Magnum"
$ws.Range("I58").Value = "Primera"
$ws.Range("J58").Value = 540
$ws.Range("K58").Value = 34000
$ws.Range("L58").Value = 35000
$ws.Range("M58").Value = 34500
$ws.Range("N58").Value = "$/malla 25 kilos"
$ws.Range("O58").Value = "Perú"
$ws.Range("P58").Value = 1380
$ws.Range("Q58").Value = 25
$ws.Range("R58").Value = "Hortaliza"
